# "matching_profiles" sheet: pandas re-exported this sheet with 24 extra
# duplicate index columns ("Unnamed: 0.1" .. "Unnamed: 0.24") in front of the
# existing headers. Recreate that by inserting 24 blank columns at B (which
# pushes the existing headers - Unnamed: 0, Applicant_Name, years_of_exp,
# Key_Skills, Linkedin_Profile, GitHub_Profile, Mail_Id - from B:H to Z:AF),
# then filling the newly inserted cells with the numbered labels and copying
# the existing bold/centered/bordered header format onto them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 24 new blank columns before the existing column B. This shifts the
# current B:H header row to Z:AF, matching the diff exactly.
$ws.Range("B1:Y1").EntireColumn.Insert()

# Give the newly inserted header cells the same look as the rest of row 1
# (bold font, thin border, centered/top aligned) by copying the format from
# the cell that used to be B1 (now at Z1) - same as Excel's own
# "Insert Copied Cells" / paste-formats workflow, and reuses the existing
# style instead of creating a near-duplicate one.
$ws.Range("Z1").Copy()
$ws.Range("B1:Y1").PasteSpecial(-4122)

# Fill the new header cells with the repeated "Unnamed: 0.N" labels, counting
# down from .24 at column B to .1 at column Y, so the original "Unnamed: 0"
# (now at Z1) is the final, un-suffixed occurrence.
$columns = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y")
$labels = @(
    "Unnamed: 0.24", "Unnamed: 0.23", "Unnamed: 0.22", "Unnamed: 0.21",
    "Unnamed: 0.20", "Unnamed: 0.19", "Unnamed: 0.18", "Unnamed: 0.17",
    "Unnamed: 0.16", "Unnamed: 0.15", "Unnamed: 0.14", "Unnamed: 0.13",
    "Unnamed: 0.12", "Unnamed: 0.11", "Unnamed: 0.10", "Unnamed: 0.9",
    "Unnamed: 0.8",  "Unnamed: 0.7",  "Unnamed: 0.6",  "Unnamed: 0.5",
    "Unnamed: 0.4",  "Unnamed: 0.3",  "Unnamed: 0.2",  "Unnamed: 0.1"
)

for ($i = 0; $i -lt $columns.Length; $i++) {
    $ws.Range("$($columns[$i])1").Value = $labels[$i]
}
